$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header columns, matching the existing "Unique" header style (L4)
$ws.Range("M4").Value = "Pattern"
$ws.Range("N4").Value = "Pattern Type"

$ws.Range("L4").Copy() | Out-Null
$ws.Range("M4:N4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("M4").Value = "Pattern"
$ws.Range("N4").Value = "Pattern Type"

# Update the selection to match the new active range
$ws.Range("M4:N4").Select() | Out-Null
